# Add two new Mac-Addresses (10 new device rows) to the
# master-reg_center_device_h sheet, rows 147-156.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$startRow = 147
$endRow = 156
$startDevice = 3000166

for ($i = 0; $i -lt ($endRow - $startRow + 1); $i++) {
    $row = $startRow + $i
    $deviceId = $startDevice + $i

    $ws.Cells.Item($row, 1).Value = 10001
    $ws.Cells.Item($row, 2).Value = $deviceId
    $ws.Cells.Item($row, 3).Value = "eng"
    $ws.Cells.Item($row, 4).Value = $true
    $ws.Cells.Item($row, 5).Value = "superadmin"
    $ws.Cells.Item($row, 6).Value = "now()"
    $ws.Cells.Item($row, 7).Value = "now()"
}

# Update the selection to match the author's final on-screen state.
$ws.Range("E155").Select()
